$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names (column B) and emails (column C) for rows 2-12.
$ws.Range("B2").Value = "Darshan"
$ws.Range("C2").Value = "pasne.d@husky.neu.edu"

$ws.Range("B3").Value = "Saman"
$ws.Range("C3").Value = "sood.s@husky.neu.edu"

$ws.Range("B4").Value = "Shail"
$ws.Range("C4").Value = "shail@ccs.neu.edu"

$ws.Range("B5").Value = "Vaibhav"
$ws.Range("C5").Value = "dave.v@husky.neu.edu"

$ws.Range("B6").Value = "John"
$ws.Range("C6").Value = "snow.j@husky.neu.edu"

$ws.Range("B7").Value = "Danny"
$ws.Range("C7").Value = "danny.d@husky.neu.edu"

$ws.Range("B8").Value = "Erica"
$ws.Range("C8").Value = "sniper.e@husky.neu.edu"

$ws.Range("B9").Value = "Flurry"
$ws.Range("C9").Value = "majin.f@husky.neu.edu"

$ws.Range("B10").Value = "Gara"
$ws.Range("C10").Value = "hawking.g@husky.neu.edu"

$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"

$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"

# Remove the mailto hyperlinks that were attached to the email column
# (the underlying cell style is preserved).
$ws.Hyperlinks.Delete()

# Move the active selection to C19 (matches the post-edit selection state).
$ws.Range("C19").Select()

$wb.Save()
